$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Needle" row (row 8) and the "Thermometer" row (originally row 15,
# which becomes row 14 once row 8 has been removed and everything below shifted up).
$ws.Rows(8).Delete()
$ws.Rows(14).Delete()

# Update the batch id text.
$ws.Range("B4").Value = "BIE24"

# Renumber the "Item id" column sequentially now that two rows were removed.
$ws.Range("A8").Value = 1
$ws.Range("A9").Value = 2
$ws.Range("A10").Value = 3
$ws.Range("A11").Value = 4
$ws.Range("A12").Value = 5
$ws.Range("A13").Value = 6
$ws.Range("A14").Value = 7

# Move the selection to B4 to match the saved view state.
$ws.Range("B4").Select() | Out-Null
